# Added New Mac-Address and Document Types
# Appends 5 new reg_center_machine_device rows (regcntr_id 10002 / machine_id
# 10032) with device_ids 3000176-3000180 to the bottom of the test-data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDeviceIds = @(3000176, 3000177, 3000178, 3000179, 3000180)

$firstNewRow = 157
for ($i = 0; $i -lt $newDeviceIds.Count; $i++) {
    $row = $firstNewRow + $i

    $ws.Cells.Item($row, 1).Value = 10002            # regcntr_id
    $ws.Cells.Item($row, 2).Value = 10032             # machine_id
    $ws.Cells.Item($row, 3).Value = $newDeviceIds[$i] # device_id
    $ws.Cells.Item($row, 4).Value = "eng"             # lang_code
    $ws.Cells.Item($row, 5).Value = $true             # is_active
    $ws.Cells.Item($row, 6).Value = "superadmin"      # cr_by
    $ws.Cells.Item($row, 7).Value = "now()"           # cr_dtimes
}

# Match the author's final cursor position/selection after the edit.
$ws.Range("E157").Select()
